$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 rows for the "ECs" sending-cluster (rows 2-4), shifting remaining rows up
$ws.Range("A2:T4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Overwrite the remaining 3 rows with the recalculated TPM-derived values

# Row 2
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.23247666666667
$ws.Range("H2").Value = 60.69743
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.4188640502130462
$ws.Range("P2").Value = 0.4188640502130463
$ws.Range("Q2").Value = 11.54504234721556
$ws.Range("R2").Value = 103.90538112494
$ws.Range("S2").Value = 0.4188640502130462
$ws.Range("T2").Value = 0.4188640502130463

# Row 3
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.23247666666667
$ws.Range("H3").Value = 60.69743
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4846943333333333
$ws.Range("N3").Value = 1.454083
$ws.Range("O3").Value = 0.3557906641356566
$ws.Range("P3").Value = 0.3557906641356566
$ws.Range("Q3").Value = 9.806566789632221
$ws.Range("R3").Value = 88.25910110668998
$ws.Range("S3").Value = 0.3557906641356566
$ws.Range("T3").Value = 0.3557906641356566

# Row 4
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.23247666666667
$ws.Range("H4").Value = 60.69743
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3069883333333334
$ws.Range("N4").Value = 0.920965
$ws.Range("O4").Value = 0.2253452856512971
$ws.Range("P4").Value = 0.2253452856512971
$ws.Range("Q4").Value = 6.211134291105556
$ws.Range("R4").Value = 55.90020861995
$ws.Range("S4").Value = 0.2253452856512971
$ws.Range("T4").Value = 0.2253452856512971
